$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 ("I0") and J1 ("IF"), copying the header style from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new I/J data columns for rows 2-16
$values = @(
    @(1, 2),
    @(7, 8),
    @(5, 5),
    @(8, 9),
    @(8, 8),
    @(7, 8),
    @(2, 3),
    @(10, 11),
    @(2, 4),
    @(7, 9),
    @(4, 4),
    @(9, 9),
    @(7, 7),
    @(3, 3),
    @(5, 5)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
